$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before the existing "Late" column (N),
# pushing Late/Heading/Outstanding one column to the right (N->O, O->P, P->Q)
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = 9.14

# Make "Repayment Schedule" the active sheet (was "Transactions"),
# with the new selection on it
$ws.Activate() | Out-Null
$ws.Range("U8").Select() | Out-Null
